$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.523.16"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "2.228.83"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'112.51"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "'294.54"
$ws.Range("E6").Value = "  +9.36%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "'43.52"
$ws.Range("E10").Value = "  -5.73%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "'54.17"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "'8.71"
$ws.Range("E13").Value = "  -4.47%  "
$ws.Range("E14").Value = "  +21.39%  "
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("D17").Value = "2.565.13"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "2.225.42"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "42.499.18"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "'7.25"
$ws.Range("E20").Value = "  +7.45%  "
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").Value = "'73.65"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").Value = "'3.37"
$ws.Range("E23").Value = "  +15.04%  "
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "'239.72"
$ws.Range("E25").Value = "  +3.81%  "
$ws.Range("D26").Value = "'8.89"
$ws.Range("E26").Value = "  -4.47%  "
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").Value = "'11.46"
$ws.Range("E28").Value = "  -6.41%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'175.29"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'37.16"
$ws.Range("E31").Value = "  -7.72%  "
$ws.Range("D32").Value = "'21.67"
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("D35").Value = "'5.69"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").Value = "'4.19"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("D41").Value = "'2.39"
$ws.Range("E41").Value = "  -6.06%  "
$ws.Range("D42").Value = "'71.40"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'12.32"
$ws.Range("E45").Value = "  -6.27%  "
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("E47").Value = "  -4.15%  "
$ws.Range("D48").Value = "'1.28"
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("D49").Value = "'8.52"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "'102.33"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("E51").Value = "  -1.21%  "
